$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 62.625
$ws.Range("I9").Value = 50
$ws.Range("K9").Value = 50
$ws.Range("M9").Value = 119

$ws.Range("H43").Value = 4966.3335
$ws.Range("I43").Value = 4966.3335
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 4966.3335
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -4897.3335
$ws.Range("N43").ClearContents()

$ws.Range("H58").Value = 117.125
$ws.Range("I58").Value = 117.125
$ws.Range("K58").Value = 351.375
$ws.Range("M58").Value = -201.375

$ws.Range("H96").Value = 47621970
$ws.Range("I96").Value = 3314.0625
$ws.Range("K96").Value = 9942.1875
$ws.Range("M96").Value = -8569.1875

$ws.Range("H107").Value = 62504704
$ws.Range("J107").Value = 6277.5
$ws.Range("L107").Value = 6277.5
$ws.Range("N107").Value = -10117.5

$ws.Range("H111").Value = 2394.1428
$ws.Range("I111").Value = 2259.5
$ws.Range("J111").Value = 2448
$ws.Range("K111").Value = 6778.5
$ws.Range("L111").Value = 7344
$ws.Range("M111").Value = -3711.5
$ws.Range("N111").Value = -13478

$ws.Range("H132").Value = 6354.0146
$ws.Range("I132").Value = 4419.9434
$ws.Range("J132").Value = 13187.733
$ws.Range("K132").Value = 13259.8302
$ws.Range("L132").Value = 39563.199
$ws.Range("M132").Value = -10729.8302
$ws.Range("N132").Value = -44623.199

$ws.Range("H135").Value = 1517.1818
$ws.Range("I135").Value = 1168.9
$ws.Range("K135").Value = 10520.1
$ws.Range("M135").Value = -7985.1

$ws.Range("H138").Value = 3513.67
$ws.Range("J138").Value = 3741.4487
$ws.Range("L138").Value = 11224.3461
$ws.Range("N138").Value = -21504.3461

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17544.646
$ws.Range("I32").Value = 12890.478
$ws.Range("K32").Value = 12890.478
$ws.Range("M32").Value = -12603.478

$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H63").Value = 2580.7144
$ws.Range("I63").Value = 2344.1667
$ws.Range("K63").Value = 2344.1667
$ws.Range("M63").Value = -1658.1667

$ws.Range("H66").Value = 2580.7144
$ws.Range("I66").Value = 2344.1667
$ws.Range("K66").Value = 11720.8335
$ws.Range("M66").Value = -8288.833500000001

$ws.Range("H122").Value = 6930
$ws.Range("J122").Value = 5482.3335
$ws.Range("L122").Value = 16447.0005
$ws.Range("N122").Value = -21347.0005

$ws.Range("H132").Value = 2675.9285
$ws.Range("I132").Value = 2461.9473
$ws.Range("J132").Value = 3127.6667
$ws.Range("K132").Value = 7385.841899999999
$ws.Range("L132").Value = 9383.000100000001
$ws.Range("M132").Value = -4855.841899999999
$ws.Range("N132").Value = -14443.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 25642304
$ws.Range("I64").Value = 47619930
$ws.Range("K64").Value = 47619930
$ws.Range("M64").Value = -47619705

$ws.Range("H67").Value = 25642304
$ws.Range("I67").Value = 47619930
$ws.Range("K67").Value = 47619930
$ws.Range("M67").Value = -47619150

$ws.Range("H117").Value = 149499
$ws.Range("J117").Value = 149499
$ws.Range("L117").Value = 149499
$ws.Range("N117").Value = -158677

$ws.Range("H134").Value = 3124.2808
$ws.Range("I134").Value = 2845.125
$ws.Range("K134").Value = 8535.375
$ws.Range("M134").Value = -6000.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4835.4863
$ws.Range("I31").Value = 3675.9167
$ws.Range("J31").Value = 5392.08
$ws.Range("K31").Value = 3675.9167
$ws.Range("L31").Value = 5392.08
$ws.Range("M31").Value = -3380.9167
$ws.Range("N31").Value = -5982.08

$ws.Range("H34").Value = 4835.4863
$ws.Range("I34").Value = 3675.9167
$ws.Range("J34").Value = 5392.08
$ws.Range("K34").Value = 3675.9167
$ws.Range("L34").Value = 5392.08
$ws.Range("M34").Value = -3473.9167
$ws.Range("N34").Value = -5796.08

$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()

$ws.Range("H58").Value = 2359.2104
$ws.Range("I58").Value = 1891.1111
$ws.Range("J58").Value = 2780.5
$ws.Range("K58").Value = 1891.1111
$ws.Range("L58").Value = 2780.5
$ws.Range("M58").Value = -1688.1111
$ws.Range("N58").Value = -3186.5

$ws.Range("H99").Value = 10877.656
$ws.Range("I99").Value = 7911.5
$ws.Range("J99").Value = 13843.8125
$ws.Range("K99").Value = 7911.5
$ws.Range("L99").Value = 13843.8125
$ws.Range("M99").Value = -6413.5
$ws.Range("N99").Value = -16839.8125

$ws.Range("H126").Value = 10877.656
$ws.Range("I126").Value = 7911.5
$ws.Range("J126").Value = 13843.8125
$ws.Range("K126").Value = 23734.5
$ws.Range("L126").Value = 41531.4375
$ws.Range("M126").Value = -21264.5
$ws.Range("N126").Value = -46471.4375

$ws.Range("H134").Value = 2754.2373
$ws.Range("I134").Value = 2448.0208
$ws.Range("J134").Value = 4090.4546
$ws.Range("K134").Value = 7344.062399999999
$ws.Range("L134").Value = 12271.3638
$ws.Range("M134").Value = -4809.062399999999
$ws.Range("N134").Value = -17341.3638

$ws.Range("H136").Value = 2359.2104
$ws.Range("I136").Value = 1891.1111
$ws.Range("J136").Value = 2780.5
$ws.Range("K136").Value = 5673.3333
$ws.Range("L136").Value = 8341.5
$ws.Range("M136").Value = -3123.3333
$ws.Range("N136").Value = -13441.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2284.4285
$ws.Range("J5").Value = 2922.25
$ws.Range("L5").Value = 8766.75
$ws.Range("N5").Value = -8990.75

$ws.Range("H61").Value = 130.41667
$ws.Range("I61").Value = 124.9
$ws.Range("J61").Value = 158
$ws.Range("K61").Value = 374.7
$ws.Range("L61").Value = 474
$ws.Range("M61").Value = -159.7
$ws.Range("N61").Value = -904

$ws.Range("H135").Value = 2284.4285
$ws.Range("J135").Value = 2922.25
$ws.Range("L135").Value = 26300.25
$ws.Range("N135").Value = -31370.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 1238.1538
$ws.Range("I17").Value = 1238.1538
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 1238.1538
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -1070.1538
$ws.Range("N17").ClearContents()

$ws.Range("H70").Value = 5549.375
$ws.Range("I70").Value = 4734.7144
$ws.Range("J70").Value = 6183
$ws.Range("K70").Value = 4734.7144
$ws.Range("L70").Value = 6183
$ws.Range("M70").Value = -4464.7144
$ws.Range("N70").Value = -6723

$ws.Range("H73").Value = 5549.375
$ws.Range("I73").Value = 4734.7144
$ws.Range("J73").Value = 6183
$ws.Range("K73").Value = 4734.7144
$ws.Range("L73").Value = 6183
$ws.Range("M73").Value = -3798.7144
$ws.Range("N73").Value = -8055

$ws.Range("H102").Value = 5706.646
$ws.Range("I102").Value = 5137.4736
$ws.Range("K102").Value = 5137.4736
$ws.Range("M102").Value = -3515.4736

$ws.Range("H113").Value = 4091.6875
$ws.Range("I113").Value = 3437.25
$ws.Range("J113").Value = 4746.125
$ws.Range("K113").Value = 3437.25
$ws.Range("L113").Value = 4746.125
$ws.Range("M113").Value = -1267.25
$ws.Range("N113").Value = -9086.125

$ws.Range("H134").Value = 166934460
$ws.Range("J134").Value = 166934460
$ws.Range("L134").Value = 500803380
$ws.Range("N134").Value = -500808450

$ws.Range("H135").Value = 105008.5
$ws.Range("J135").Value = 105008.5
$ws.Range("L135").Value = 105008.5
$ws.Range("N135").Value = -115148.5

$ws.Range("H136").Value = 98550.336
$ws.Range("J136").Value = 98550.336
$ws.Range("L136").Value = 295651.008
$ws.Range("N136").Value = -300751.008

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6164.75
$ws.Range("I40").Value = 4533.647
$ws.Range("J40").Value = 10126
$ws.Range("K40").Value = 4533.647
$ws.Range("L40").Value = 10126
$ws.Range("M40").Value = -4397.647
$ws.Range("N40").Value = -10398

$ws.Range("H46").Value = 1379.8572
$ws.Range("I46").Value = 943.3333
$ws.Range("J46").Value = 1498.909
$ws.Range("K46").Value = 943.3333
$ws.Range("L46").Value = 1498.909
$ws.Range("M46").Value = -755.3333
$ws.Range("N46").Value = -1874.909

$ws.Range("H55").Value = 732.5789
$ws.Range("J55").Value = 869.4
$ws.Range("L55").Value = 869.4
$ws.Range("N55").Value = -1215.4

$ws.Range("H136").Value = 12308.348
$ws.Range("I136").Value = 14524.75
$ws.Range("J136").Value = 11126.267
$ws.Range("K136").Value = 43574.25
$ws.Range("L136").Value = 33378.801
$ws.Range("M136").Value = -41024.25
$ws.Range("N136").Value = -38478.801

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 15922.625
$ws.Range("J41").Value = 15983
$ws.Range("L41").Value = 15983
$ws.Range("N41").Value = -16763

$ws.Range("H107").Value = 735.0513
$ws.Range("I107").Value = 654.73334
$ws.Range("K107").Value = 1964.20002
$ws.Range("M107").Value = -44.20001999999999

$ws.Range("H122").Value = 8214.75
$ws.Range("I122").Value = 6148.852
$ws.Range("K122").Value = 18446.556
$ws.Range("M122").Value = -15996.556

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
